$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New actual data point for Mar 30 lands in row 25 (column I), replacing the
# previous forecast formula with a hardcoded actual value. Pick up the
# "actual data" formatting (matches I21:I24) instead of the forecast shading.
$ws.Range("I24").Copy()
$ws.Range("I25").PasteSpecial(-4122)
$ws.Range("I25").Value = 163788

# The next forecast row (I26) now anchors directly off the new actual value
# with an updated rolling-average window.
$ws.Range("I26").Formula = "=I25*(1+AVERAGE(M24:M25))"

# Reflect the new active selection left by the edit.
$ws.Range("I27").Select()

$wb.Application.Calculate()
